# Updated symbol list on Fri Jan  6 21:54:01 UTC 2023 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) text values for the crypto ranking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "259.03"
Set-TextValue "E2" "0.61%"
Set-TextValue "D3" "27.01"
Set-TextValue "E3" "0.18%"
Set-TextValue "D4" "4.702"
Set-TextValue "E4" "0.42%"
Set-TextValue "D5" "0.06032"
Set-TextValue "E5" "2.51%"
Set-TextValue "E6" "0.43%"
Set-TextValue "D7" "0.8600"
Set-TextValue "E7" "0.26%"
Set-TextValue "D8" "0.9237"
Set-TextValue "E8" "-3.14%"
Set-TextValue "D9" "0.1396"
Set-TextValue "E9" "-0.89%"
Set-TextValue "D10" "0.05161"
Set-TextValue "E10" "30.92%"
Set-TextValue "D11" "0.07074"
Set-TextValue "E11" "-0.15%"
Set-TextValue "E12" "-2.95%"
Set-TextValue "D13" "0.09128"
Set-TextValue "E13" "-0.30%"
Set-TextValue "D14" "0.001529"
Set-TextValue "E14" "-1.51%"
Set-TextValue "D15" "0.0006042"
Set-TextValue "E15" "-94.25%"
Set-TextValue "D16" "0.006094"
Set-TextValue "E16" "-1.67%"
Set-TextValue "D17" "3.466"
Set-TextValue "E17" "-1.35%"
Set-TextValue "D18" "3.172"
Set-TextValue "E18" "-0.95%"
Set-TextValue "D19" "2.166"
Set-TextValue "E19" "-2.77%"
Set-TextValue "E20" "0.40%"
Set-TextValue "E21" "0.37%"
Set-TextValue "D23" "0.04226"
Set-TextValue "E23" "-0.11%"
Set-TextValue "D24" "0.001218"
Set-TextValue "E24" "-0.30%"
Set-TextValue "D25" "0.004037"
Set-TextValue "E26" "0.02%"
Set-TextValue "D27" "0.0001523"
Set-TextValue "E27" "-21.36%"
Set-TextValue "D40" "0.03849"
Set-TextValue "D41" "0.1114"
Set-TextValue "E41" "1.12%"
Set-TextValue "D42" "0.004021"
Set-TextValue "E42" "-35.45%"
Set-TextValue "D43" "0.01529"
Set-TextValue "E43" "33.53%"
Set-TextValue "E44" "0.02%"
Set-TextValue "D45" "0.00005101"
Set-TextValue "E45" "-6.40%"
Set-TextValue "E46" "0.02%"
Set-TextValue "D47" "0.1353"
Set-TextValue "E47" "-24.60%"
Set-TextValue "D48" "0.05454"
Set-TextValue "E48" "-22.08%"
Set-TextValue "E49" "0.02%"
Set-TextValue "E50" "0.02%"
